$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 5.690303
$ws.Range("H2").Value = 17.070909
$ws.Range("I2").Value = 0.3759010823723209
$ws.Range("J2").Value = 0.3759010823723208
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.971099
$ws.Range("N2").Value = 2.913297
$ws.Range("O2").Value = 0.007882574716876797
$ws.Range("P2").Value = 0.007882574716876797
$ws.Range("Q2").Value = 5.525847552997
$ws.Range("R2").Value = 49.732627976973
$ws.Range("S2").Value = 0.002963068367954679
$ws.Range("T2").Value = 0.002963068367954679
# Row 3
$ws.Range("G3").Value = 5.690303
$ws.Range("H3").Value = 17.070909
$ws.Range("I3").Value = 0.3759010823723209
$ws.Range("J3").Value = 0.3759010823723208
$ws.Range("O3").Value = 0.6966643430097871
$ws.Range("P3").Value = 0.696664343009787
$ws.Range("Q3").Value = 488.376081845272
$ws.Range("R3").Value = 4395.384736607448
$ws.Range("S3").Value = 0.2618768805875808
$ws.Range("T3").Value = 0.2618768805875807
# Row 4
$ws.Range("G4").Value = 5.690303
$ws.Range("H4").Value = 17.070909
$ws.Range("I4").Value = 0.3759010823723209
$ws.Range("J4").Value = 0.3759010823723208
$ws.Range("M4").Value = 36.24916566666667
$ws.Range("N4").Value = 108.747497
$ws.Range("O4").Value = 0.294240604502677
$ws.Range("P4").Value = 0.294240604502677
$ws.Range("Q4").Value = 206.2687361405304
$ws.Range("R4").Value = 1856.418625264773
$ws.Range("S4").Value = 0.1106053617104423
$ws.Range("T4").Value = 0.1106053617104423
# Row 5
$ws.Range("G5").Value = 5.690303
$ws.Range("H5").Value = 17.070909
$ws.Range("I5").Value = 0.3759010823723209
$ws.Range("J5").Value = 0.3759010823723208
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.149372
$ws.Range("N5").Value = 0.448116
$ws.Range("O5").Value = 0.001212477770659141
$ws.Range("P5").Value = 0.001212477770659141
$ws.Range("Q5").Value = 0.8499719397159999
$ws.Range("R5").Value = 7.649747457444
$ws.Range("S5").Value = 0.0004557717063431497
$ws.Range("T5").Value = 0.0004557717063431495
# Row 6
$ws.Range("I6").Value = 0.2572281411562661
$ws.Range("J6").Value = 0.2572281411562661
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.971099
$ws.Range("N6").Value = 2.913297
$ws.Range("O6").Value = 0.007882574716876797
$ws.Range("P6").Value = 0.007882574716876797
$ws.Range("Q6").Value = 3.781323228440333
$ws.Range("R6").Value = 34.031909055963
$ws.Range("S6").Value = 0.002027620041947599
$ws.Range("T6").Value = 0.002027620041947599
# Row 7
$ws.Range("I7").Value = 0.2572281411562661
$ws.Range("J7").Value = 0.2572281411562661
$ws.Range("O7").Value = 0.6966643430097871
$ws.Range("P7").Value = 0.696664343009787
$ws.Range("S7").Value = 0.1792016739622589
$ws.Range("T7").Value = 0.1792016739622589
# Row 8
$ws.Range("I8").Value = 0.2572281411562661
$ws.Range("J8").Value = 0.2572281411562661
$ws.Range("M8").Value = 36.24916566666667
$ws.Range("N8").Value = 108.747497
$ws.Range("O8").Value = 0.294240604502677
$ws.Range("P8").Value = 0.294240604502677
$ws.Range("Q8").Value = 141.1491641397514
$ws.Range("R8").Value = 1270.342477257763
$ws.Range("S8").Value = 0.07568696374891967
$ws.Range("T8").Value = 0.07568696374891966
# Row 9
$ws.Range("I9").Value = 0.2572281411562661
$ws.Range("J9").Value = 0.2572281411562661
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.149372
$ws.Range("N9").Value = 0.448116
$ws.Range("O9").Value = 0.001212477770659141
$ws.Range("P9").Value = 0.001212477770659141
$ws.Range("Q9").Value = 0.5816336061293332
$ws.Range("R9").Value = 5.234702455163999
$ws.Range("S9").Value = 0.0003118834031399442
$ws.Range("T9").Value = 0.0003118834031399442
# Row 10
$ws.Range("G10").Value = 4.648693333333333
$ws.Range("H10").Value = 13.94608
$ws.Range("I10").Value = 0.3070924088957991
$ws.Range("J10").Value = 0.307092408895799
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.971099
$ws.Range("N10").Value = 2.913297
$ws.Range("O10").Value = 0.007882574716876797
$ws.Range("P10").Value = 0.007882574716876797
$ws.Range("Q10").Value = 4.514341447306666
$ws.Range("R10").Value = 40.62907302575999
$ws.Range("S10").Value = 0.002420678858106817
$ws.Range("T10").Value = 0.002420678858106817
# Row 11
$ws.Range("G11").Value = 4.648693333333333
$ws.Range("H11").Value = 13.94608
$ws.Range("I11").Value = 0.3070924088957991
$ws.Range("J11").Value = 0.307092408895799
$ws.Range("O11").Value = 0.6966643430097871
$ws.Range("P11").Value = 0.696664343009787
$ws.Range("Q11").Value = 398.9788655953067
$ws.Range("R11").Value = 3590.80979035776
$ws.Range("S11").Value = 0.2139403312866848
$ws.Range("T11").Value = 0.2139403312866847
# Row 12
$ws.Range("G12").Value = 4.648693333333333
$ws.Range("H12").Value = 13.94608
$ws.Range("I12").Value = 0.3070924088957991
$ws.Range("J12").Value = 0.307092408895799
$ws.Range("M12").Value = 36.24916566666667
$ws.Range("N12").Value = 108.747497
$ws.Range("O12").Value = 0.294240604502677
$ws.Range("P12").Value = 0.294240604502677
$ws.Range("Q12").Value = 168.5112547735289
$ws.Range("R12").Value = 1516.60129296176
$ws.Range("S12").Value = 0.09035905603168319
$ws.Range("T12").Value = 0.09035905603168316
# Row 13
$ws.Range("G13").Value = 4.648693333333333
$ws.Range("H13").Value = 13.94608
$ws.Range("I13").Value = 0.3070924088957991
$ws.Range("J13").Value = 0.307092408895799
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 0.6666666666666666
$ws.Range("M13").Value = 0.149372
$ws.Range("N13").Value = 0.448116
$ws.Range("O13").Value = 0.001212477770659141
$ws.Range("P13").Value = 0.001212477770659141
$ws.Range("Q13").Value = 0.6943846205866665
$ws.Range("R13").Value = 6.249461585279999
$ws.Range("S13").Value = 0.0003723427193243237
$ws.Range("T13").Value = 0.0003723427193243237
# Row 14
$ws.Range("G14").Value = 0.9049109999999999
$ws.Range("H14").Value = 2.714733
$ws.Range("I14").Value = 0.05977836757561403
$ws.Range("J14").Value = 0.05977836757561403
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 0.971099
$ws.Range("N14").Value = 2.913297
$ws.Range("O14").Value = 0.007882574716876797
$ws.Range("P14").Value = 0.007882574716876797
$ws.Range("Q14").Value = 0.8787581671889999
$ws.Range("R14").Value = 7.908823504701
$ws.Range("S14").Value = 0.0004712074488677028
$ws.Range("T14").Value = 0.0004712074488677028
# Row 15
$ws.Range("G15").Value = 0.9049109999999999
$ws.Range("H15").Value = 2.714733
$ws.Range("I15").Value = 0.05977836757561403
$ws.Range("J15").Value = 0.05977836757561403
$ws.Range("O15").Value = 0.6966643430097871
$ws.Range("P15").Value = 0.696664343009787
$ws.Range("Q15").Value = 77.66491320386399
$ws.Range("R15").Value = 698.9842188347759
$ws.Range("S15").Value = 0.04164545717326271
$ws.Range("T15").Value = 0.0416454571732627
# Row 16
$ws.Range("G16").Value = 0.9049109999999999
$ws.Range("H16").Value = 2.714733
$ws.Range("I16").Value = 0.05977836757561403
$ws.Range("J16").Value = 0.05977836757561403
$ws.Range("M16").Value = 36.24916566666667
$ws.Range("N16").Value = 108.747497
$ws.Range("O16").Value = 0.294240604502677
$ws.Range("P16").Value = 0.294240604502677
$ws.Range("Q16").Value = 32.802268752589
$ws.Range("R16").Value = 295.220418773301
$ws.Range("S16").Value = 0.0175892230116319
$ws.Range("T16").Value = 0.01758922301163189
# Row 17
$ws.Range("G17").Value = 0.9049109999999999
$ws.Range("H17").Value = 2.714733
$ws.Range("I17").Value = 0.05977836757561403
$ws.Range("J17").Value = 0.05977836757561403
$ws.Range("K17").Value = 2
$ws.Range("L17").Value = 0.6666666666666666
$ws.Range("M17").Value = 0.149372
$ws.Range("N17").Value = 0.448116
$ws.Range("O17").Value = 0.001212477770659141
$ws.Range("P17").Value = 0.001212477770659141
$ws.Range("Q17").Value = 0.135168365892
$ws.Range("R17").Value = 1.216515293028
$ws.Range("S17").Value = 0.00007247994185172315
$ws.Range("T17").Value = 0.00007247994185172315
